$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRange, $text) {
    $cellRange.NumberFormat = "@"
    $cellRange.Value = $text
    $cellRange.ClearFormats()
}

Set-TextValue $ws.Range("D2") "27.181.04"
$ws.Range("E2").Value = "  +0.91%  "
Set-TextValue $ws.Range("D3") "1.901.90"
$ws.Range("E3").Value = "  +1.37%  "
$ws.Range("E4").Value = "  +0.23%  "
Set-TextValue $ws.Range("D5") "306.70"
$ws.Range("E5").Value = "  +0.11%  "
$ws.Range("E6").Value = "  +0.16%  "
Set-TextValue $ws.Range("D7") "0.5229"
$ws.Range("E7").Value = "  +1.59%  "
$ws.Range("E8").Value = "  +0.69%  "
Set-TextValue $ws.Range("D9") "0.07240"
$ws.Range("E9").Value = "  +0.73%  "
Set-TextValue $ws.Range("D10") "21.18"
$ws.Range("E10").Value = "  +2.14%  "
Set-TextValue $ws.Range("D11") "0.8982"
$ws.Range("E11").Value = "  -0.11%  "
Set-TextValue $ws.Range("D12") "0.08434"
$ws.Range("E12").Value = "  +11.69%  "
Set-TextValue $ws.Range("D13") "1.901.97"
$ws.Range("E13").Value = "  +1.47%  "
Set-TextValue $ws.Range("D14") "94.63"
$ws.Range("E14").Value = "  -0.57%  "
Set-TextValue $ws.Range("D15") "5.264"
$ws.Range("E15").Value = "  +0.18%  "
$ws.Range("E16").Value = "  +0.22%  "
Set-TextValue $ws.Range("D17") "0.000008582"
$ws.Range("E17").Value = "  +1.12%  "
$ws.Range("E18").Value = "  +1.66%  "
$ws.Range("E19").Value = "  +0.14%  "
Set-TextValue $ws.Range("D20") "27.224.47"
$ws.Range("E20").Value = "  +0.91%  "
$ws.Range("E21").Value = "  +0.57%  "
Set-TextValue $ws.Range("D22") "2.146.20"
$ws.Range("E22").Value = "  +0.26%  "
$ws.Range("E23").Value = "  +1.62%  "
Set-TextValue $ws.Range("D24") "6.417"
$ws.Range("E24").Value = "  +0.00%  "
Set-TextValue $ws.Range("D25") "146.65"
$ws.Range("E25").Value = "  +0.46%  "
Set-TextValue $ws.Range("D26") "2.276"
$ws.Range("E26").Value = "  +7.54%  "
Set-TextValue $ws.Range("D27") "1.751"
$ws.Range("E27").Value = "  -1.72%  "
$ws.Range("E28").Value = "  +0.85%  "
Set-TextValue $ws.Range("D29") "114.82"
$ws.Range("E29").Value = "  -0.21%  "
Set-TextValue $ws.Range("D30") "4.921"
$ws.Range("E30").Value = "  -0.30%  "
Set-TextValue $ws.Range("D31") "4.782"
$ws.Range("E31").Value = "  +0.16%  "
$ws.Range("E32").Value = "  +0.31%  "
Set-TextValue $ws.Range("D33") "0.8117"
$ws.Range("E33").Value = "  +7.99%  "
Set-TextValue $ws.Range("D34") "0.05060"
$ws.Range("E34").Value = "  +0.62%  "
$ws.Range("E35").Value = "  +5.33%  "
Set-TextValue $ws.Range("D36") "2.960"
$ws.Range("E36").Value = "  -0.95%  "
Set-TextValue $ws.Range("D37") "3.370"
$ws.Range("E37").Value = "  +3.66%  "
Set-TextValue $ws.Range("D38") "2.559"
$ws.Range("E38").Value = "  +2.81%  "
$ws.Range("E39").Value = "  +1.89%  "
Set-TextValue $ws.Range("D40") "0.01978"
$ws.Range("E40").Value = "  -0.54%  "
Set-TextValue $ws.Range("D41") "1.072"
$ws.Range("E41").Value = "  +0.14%  "
$ws.Range("B42").Value = "Aptos"
$ws.Range("C42").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue $ws.Range("D42") "8.951"
$ws.Range("E42").Value = "  +3.05%  "
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue $ws.Range("D43") "6.614"
$ws.Range("E43").Value = "  +0.64%  "
Set-TextValue $ws.Range("D44") "118.32"
$ws.Range("E44").Value = "  +2.35%  "
Set-TextValue $ws.Range("D45") "0.1508"
$ws.Range("E45").Value = "  +0.61%  "
Set-TextValue $ws.Range("D46") "0.4824"
$ws.Range("E46").Value = "  +0.93%  "
$ws.Range("E47").Value = "  +0.24%  "
Set-TextValue $ws.Range("D48") "10.12"
$ws.Range("E48").Value = "  +0.06%  "
Set-TextValue $ws.Range("D49") "1.611"
$ws.Range("E49").Value = "  +2.77%  "
Set-TextValue $ws.Range("D50") "37.40"
$ws.Range("E50").Value = "  +1.08%  "
Set-TextValue $ws.Range("D51") "63.53"
$ws.Range("E51").Value = "  +0.43%  "
